$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style as the other header cells (B1:H1) to the new headers
$ws.Range("B1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-26
$values = @{
    2  = @(5, 6)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(6, 7)
    8  = @(8, 8)
    9  = @(9, 9)
    10 = @(5, 6)
    11 = @(6, 6)
    12 = @(8, 8)
    13 = @(6, 6)
    14 = @(6, 8)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(6, 7)
    18 = @(4, 5)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(6, 6)
    22 = @(4, 4)
    23 = @(9, 9)
    24 = @(9, 9)
    25 = @(8, 9)
    26 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
